$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1212, shifting existing rows 1212-1315 down to 1213-1316
$ws.Rows.Item(1212).Insert()

# Populate the newly inserted row 1212 with the new record's data
$ws.Range("A1212").Value = 6
$ws.Range("B1212").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1212").Value = "Metropolitana"
$ws.Range("D1212").Value = 45106
$ws.Range("E1212").Value = 13
$ws.Range("F1212").Value = 100112031
$ws.Range("G1212").Value = "Poroto verde"
$ws.Range("H1212").Value = "Magnum"
$ws.Range("I1212").Value = "Primera"
$ws.Range("J1212").Value = 500
$ws.Range("K1212").Value = 15000
$ws.Range("L1212").Value = 17000
$ws.Range("M1212").Value = 16080
$ws.Range("N1212").Value = "$/malla 25 kilos"
$ws.Range("O1212").Value = "Perú"
$ws.Range("P1212").Value = 643
$ws.Range("Q1212").Value = 25
$ws.Range("R1212").Value = "Hortaliza"
